# Swap the presentation's theme color scheme from the custom "Integral"
# palette over to the stock "Office" palette (dk1/lt1 unchanged; dk2, lt2,
# accent1-6, hlink and folHlink are replaced).
#
# PowerPoint's ThemeColor.RGB (and the legacy RGB() macro) packs colors as
# 0xBBGGRR, so each target 0xRRGGBB hex code below is byte-swapped before
# being assigned.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> target color (standard RRGGBB, matching the "Office" theme)
$targets = @{
    1  = 0x000000  # dk1
    2  = 0xFFFFFF  # lt1
    3  = 0x44546A  # dk2
    4  = 0xE7E6E6  # lt2
    5  = 0x5B9BD5  # accent1
    6  = 0xED7D31  # accent2
    7  = 0xA5A5A5  # accent3
    8  = 0xFFC000  # accent4
    9  = 0x4472C4  # accent5
    10 = 0x70AD47  # accent6
    11 = 0x0563C1  # hlink
    12 = 0x954F72  # folHlink
}

for ($i = 1; $i -le $tcs.Count; $i++) {
    $rgb = $targets[$i]
    $r = ($rgb -shr 16) -band 0xFF
    $g = ($rgb -shr 8) -band 0xFF
    $b = $rgb -band 0xFF
    $bgr = ($b -shl 16) -bor ($g -shl 8) -bor $r
    $tcs.Colors($i).RGB = $bgr
}
